$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M2").Value = 1105.44
$ws1.Range("D6").Value = 1900.8
$ws1.Range("M6").Value = 3108.61
$ws1.Range("D58").Value = "2 de 56"
$ws1.Range("M58").Value = "9 de 56"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F2").Value = 1105.44
$ws2.Range("F6").Value = 5009.41
$ws2.Range("F58").Value = 23154.28

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 2816.64
$ws3.Range("E3").Value = 2175.543200000001
$ws3.Range("F3").Value = 0.564210063444787

$ws3.Range("D16").Value = 16405.97
$ws3.Range("E16").Value = 23984.2
$ws3.Range("F16").Value = 0.4061871985188476

$ws3.Range("D19").Value = 23154.28
$ws3.Range("E19").Value = 32255.42560036207
$ws3.Range("F19").Value = 0.4178740844970073
